$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.426.19"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.608.56"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "1.833.99"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "1.609.76"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "234.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.08%  "
$ws.Range("D18").Value = "26.413.56"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +4.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0495"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "1.491.11"
$ws.Range("E32").Value = "  +5.48%  "
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.563"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.926"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").Value = "1.746.58"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.760"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.11%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0501"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.82%  "
